$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "92÷6=15, 2"
$t.Cell(1, 2).Range.Text = "83÷6=13, 5"
$t.Cell(1, 3).Range.Text = "40÷4=10, 0"
$t.Cell(1, 4).Range.Text = "76÷6=12, 4"
$t.Cell(1, 5).Range.Text = "58÷7=8, 2"
$t.Cell(5, 1).Range.Text = "22÷4=5, 2"
$t.Cell(5, 2).Range.Text = "77÷6=12, 5"
$t.Cell(5, 3).Range.Text = "85÷4=21, 1"
$t.Cell(5, 4).Range.Text = "99÷2=49, 1"
$t.Cell(5, 5).Range.Text = "90÷6=15, 0"
$t.Cell(9, 1).Range.Text = "48÷4=12, 0"
$t.Cell(9, 2).Range.Text = "97÷2=48, 1"
$t.Cell(9, 3).Range.Text = "74÷2=37, 0"
$t.Cell(9, 4).Range.Text = "75÷6=12, 3"
$t.Cell(9, 5).Range.Text = "20÷3=6, 2"
$t.Cell(13, 1).Range.Text = "97÷6=16, 1"
$t.Cell(13, 2).Range.Text = "22÷3=7, 1"
$t.Cell(13, 3).Range.Text = "35÷5=7, 0"
$t.Cell(13, 4).Range.Text = "98÷5=19, 3"
$t.Cell(13, 5).Range.Text = "99÷5=19, 4"
$t.Cell(17, 1).Range.Text = "82÷3=27, 1"
$t.Cell(17, 2).Range.Text = "23÷9=2, 5"
$t.Cell(17, 3).Range.Text = "39÷2=19, 1"
$t.Cell(17, 4).Range.Text = "30÷7=4, 2"
$t.Cell(17, 5).Range.Text = "91÷6=15, 1"
